$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.400.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +2.21%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.153.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.81%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'536.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +3.04%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'139.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +2.95%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.06%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.514"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +8.95%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'7.34"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.61%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.109"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +3.50%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.420"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +5.36%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.139"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +2.25%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'3.698.38"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +2.90%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'25.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +2.57%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  +6.40%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'58.470.35"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +2.31%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.162.66"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +2.97%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'6.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +6.24%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'13.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +5.05%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'8.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +6.01%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'373.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +7.37%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +1.48%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.996"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.22%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'70.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.83%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +3.44%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +1.43%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.01%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'8.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +14.12%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.0₃0866"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +3.17%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +3.21%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'6.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +5.61%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'21.90"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +4.71%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'5.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +7.96%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.17"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +5.27%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'160.68"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.57%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +4.37%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'1.38"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +13.92%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'25.26"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.13%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'2.642.09"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +9.77%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +6.46%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +4.16%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'4.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +4.72%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +5.97%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.708"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +2.85%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.0282"
$ws.Range("D45").Style = "Normal"
$ws.Range("E46").Value = "'  +0.04%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'3.199.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +2.97%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +10.75%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'6.20"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +4.31%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.977"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +5.05%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'20.21"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +4.85%  "
$ws.Range("E51").Style = "Normal"
